# Auto-generated COM-interop script to apply scheduled price-refresh update
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 92.53846
$ws.Range("I6").Value = 96.083336
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 288.250008
$ws.Range("L6").Value = 150
$ws.Range("M6").Value = -176.250008
$ws.Range("N6").Value = -374
$ws.Range("H17").Value = 1259.2069
$ws.Range("J17").Value = 1259.2069
$ws.Range("L17").Value = 3777.620699999999
$ws.Range("N17").Value = -4113.620699999999
$ws.Range("H74").Value = 6544.1333
$ws.Range("I74").Value = 6212
$ws.Range("J74").Value = 6765.5557
$ws.Range("K74").Value = 6212
$ws.Range("L74").Value = 6765.5557
$ws.Range("M74").Value = -5276
$ws.Range("N74").Value = -8637.555700000001
$ws.Range("H77").Value = 6544.1333
$ws.Range("I77").Value = 6212
$ws.Range("J77").Value = 6765.5557
$ws.Range("K77").Value = 31060
$ws.Range("L77").Value = 33827.7785
$ws.Range("M77").Value = -26380
$ws.Range("N77").Value = -43187.7785
$ws.Range("H125").Value = 1931.4615
$ws.Range("I125").Value = 1825.7142
$ws.Range("K125").Value = 16431.4278
$ws.Range("M125").Value = -13971.4278
$ws.Range("H132").Value = 3849.6345
$ws.Range("I132").Value = 3681.7346
$ws.Range("K132").Value = 11045.2038
$ws.Range("M132").Value = -8515.203799999999
$ws.Range("H138").Value = 2935.125
$ws.Range("J138").Value = 5586.6
$ws.Range("L138").Value = 16759.8
$ws.Range("N138").Value = -27039.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 168.5
$ws.Range("I5").Value = 123.71429
$ws.Range("J5").Value = 273
$ws.Range("K5").Value = 123.71429
$ws.Range("L5").Value = 273
$ws.Range("M5").Value = -11.71429000000001
$ws.Range("N5").Value = -497
$ws.Range("H74").Value = 2978.1365
$ws.Range("I74").Value = 4596.5
$ws.Range("K74").Value = 4596.5
$ws.Range("M74").Value = -3722.5
$ws.Range("H77").Value = 2978.1365
$ws.Range("I77").Value = 4596.5
$ws.Range("K77").Value = 22982.5
$ws.Range("M77").Value = -18614.5
$ws.Range("H97").Value = 3043.8
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H132").Value = 1563.7073
$ws.Range("I132").Value = 1436.3889
$ws.Range("K132").Value = 4309.1667
$ws.Range("M132").Value = -1779.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 168.5
$ws.Range("I4").Value = 123.71429
$ws.Range("J4").Value = 273
$ws.Range("K4").Value = 123.71429
$ws.Range("L4").Value = 273
$ws.Range("M4").Value = -8.714290000000005
$ws.Range("N4").Value = -503
$ws.Range("H94").Value = 5926.7837
$ws.Range("I94").Value = 3696.2964
$ws.Range("K94").Value = 3696.2964
$ws.Range("M94").Value = -3245.2964
$ws.Range("H134").Value = 6071.959
$ws.Range("I134").Value = 4569.528
$ws.Range("J134").Value = 10232.538
$ws.Range("K134").Value = 13708.584
$ws.Range("L134").Value = 30697.614
$ws.Range("M134").Value = -11173.584
$ws.Range("N134").Value = -35767.614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 111.75
$ws.Range("I7").Value = 159.88889
$ws.Range("J7").Value = 82.86667
$ws.Range("K7").Value = 159.88889
$ws.Range("L7").Value = 82.86667
$ws.Range("M7").Value = -46.88889
$ws.Range("N7").Value = -308.86667
$ws.Range("H134").Value = 1251.6177
$ws.Range("I134").Value = 1262.303
$ws.Range("K134").Value = 3786.909000000001
$ws.Range("M134").Value = -1251.909000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1500
$ws.Range("J23").Value = 1500
$ws.Range("L23").Value = 4500
$ws.Range("N23").Value = -4970
$ws.Range("H122").Value = 1734.8334
$ws.Range("J122").Value = 4725
$ws.Range("L122").Value = 42525
$ws.Range("N122").Value = -47425
$ws.Range("H125").Value = 13636.272
$ws.Range("I125").Value = 9999
$ws.Range("J125").Value = 14000
$ws.Range("K125").Value = 29997
$ws.Range("L125").Value = 42000
$ws.Range("M125").Value = -25077
$ws.Range("N125").Value = -51840
$ws.Range("H131").Value = 16462577
$ws.Range("I131").Value = 10102258
$ws.Range("J131").Value = 20835296
$ws.Range("K131").Value = 30306774
$ws.Range("L131").Value = 62505888
$ws.Range("M131").Value = -30301734
$ws.Range("N131").Value = -62515968
$ws.Range("H137").Value = 8387.333000000001
$ws.Range("I137").Value = 8011.5
$ws.Range("K137").Value = 24034.5
$ws.Range("M137").Value = -18934.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H103").Value = 65000
$ws.Range("J103").Value = 65000
$ws.Range("L103").Value = 65000
$ws.Range("N103").Value = -67344
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37080
$ws.Range("H122").Value = 3969.8096
$ws.Range("I122").Value = 3628.4
$ws.Range("K122").Value = 10885.2
$ws.Range("M122").Value = -8435.200000000001
$ws.Range("H126").Value = 6554.923
$ws.Range("I126").Value = 5558.125
$ws.Range("J126").Value = 8149.8
$ws.Range("K126").Value = 16674.375
$ws.Range("L126").Value = 24449.4
$ws.Range("M126").Value = -14204.375
$ws.Range("N126").Value = -29389.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H61").Value = 2310.7
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 2384.0557
$ws.Range("I68").Value = 2402.5715
$ws.Range("J68").Value = 2319.25
$ws.Range("K68").Value = 2402.5715
$ws.Range("L68").Value = 2319.25
$ws.Range("M68").Value = -1653.5715
$ws.Range("N68").Value = -3817.25
$ws.Range("H71").Value = 2384.0557
$ws.Range("I71").Value = 2402.5715
$ws.Range("J71").Value = 2319.25
$ws.Range("K71").Value = 12012.8575
$ws.Range("L71").Value = 11596.25
$ws.Range("M71").Value = -8268.8575
$ws.Range("N71").Value = -19084.25
$ws.Range("H113").Value = 2310.7
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 6755.875
$ws.Range("I122").Value = 6897
$ws.Range("K122").Value = 20691
$ws.Range("M122").Value = -18241
$ws.Range("H132").Value = 2340.7
$ws.Range("I132").Value = 2109.7273
$ws.Range("K132").Value = 6329.1819
$ws.Range("M132").Value = -3799.1819
$ws.Range("H136").Value = 2336.3635
$ws.Range("I136").Value = 1932.3
$ws.Range("J136").Value = 2673.0833
$ws.Range("K136").Value = 5796.9
$ws.Range("L136").Value = 8019.249899999999
$ws.Range("M136").Value = -3246.9
$ws.Range("N136").Value = -13119.2499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H81").Value = 5759.5
$ws.Range("I81").Value = 6291.6665
$ws.Range("K81").Value = 12583.333
$ws.Range("M81").Value = -11522.333
$ws.Range("H84").Value = 5759.5
$ws.Range("I84").Value = 6291.6665
$ws.Range("K84").Value = 62916.665
$ws.Range("M84").Value = -57612.665
$ws.Range("H107").Value = 1579.1428
$ws.Range("I107").Value = 1342.3334
$ws.Range("K107").Value = 4027.0002
$ws.Range("M107").Value = -2107.0002
